# Weekly data refresh: a new week's price observation is added at the top
# of the data table (row 2), pushing all existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (2..45) down to (3..46)
$ws.Rows.Item(2).Insert()

# The Insert() above copies formatting down from the old row 2 (now row 3),
# which picks up styling that doesn't belong on a fresh data row -> clear it.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with this week's record
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44860
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100102
$ws.Range("H2").Value = "Cítricos"
$ws.Range("I2").Value = 100102006
$ws.Range("J2").Value = "Pomelo"
$ws.Range("K2").Value = "Start Ruby"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = 180000
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("Q2").Value = "`$/bins (350 kilos)"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 514
$ws.Range("T2").Value = 350
